$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) cells: "<Name>_old" -> "<Name>_FV2404",
#        "<Name>_new" -> "<Name>_FV2410". The "diff" header (column K) stays the same. ---
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Freeze the header row (split/freeze after row 1, pane stays on row 2). ---
$ws.Range("A2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $true | Out-Null

# --- 3. Turn the used range into an actual Excel Table ("Table1") so the header
#        row gets a filter dropdown + structured reference support. ---
$rng = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
